$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original row heights so multi-line values don't trigger an
# unwanted auto-fit resize of the header/data rows.
$origRow1Height = $ws.Rows(1).RowHeight
$origRow2Height = $ws.Rows(2).RowHeight

# --- Extend formatting of the new header cells (row 1) to match existing headers ---
$ws.Range("C1").Copy($ws.Range("E1"))
$ws.Range("C1").Copy($ws.Range("F1"))
$ws.Range("C1").Copy($ws.Range("G1"))
$ws.Range("C1").Copy($ws.Range("H1"))

# --- Extend formatting of the new data cells (row 2) to match existing data cells ---
$ws.Range("D2").Copy($ws.Range("E2"))
$ws.Range("D2").Copy($ws.Range("F2"))
$ws.Range("D2").Copy($ws.Range("G2"))
$ws.Range("D2").Copy($ws.Range("H2"))

# --- Header row (row 1) values ---
# C1 used to be "Groups" / D1 used to be "Business name".
# New layout inserts Country/Organization/Age/Gender/Sectors before Groups (now in H1).
$ws.Range("C1").Value = "Country"
$ws.Range("D1").Value = "Organization"
$ws.Range("E1").Value = "Age"
$ws.Range("F1").Value = "Gender"
$ws.Range("G1").Value = "Sectors"
$ws.Range("H1").Value = "Groups"

# --- Data row (row 2) values ---
$ws.Range("A2").Value = "Roberto Greco"
$ws.Range("C2").Value = "Italy"
$ws.Range("D2").Value = "Square"
$ws.Range("E2").Value = "31-50"
$ws.Range("F2").Value = "m"
$ws.Range("G2").Value = "other`nIT"
$ws.Range("H2").Value = "other`ninformation"

# Restore original row heights (the multi-line G2/H2 text can trigger an
# automatic row-height resize which the source workbook does not have).
$ws.Rows(1).RowHeight = $origRow1Height
$ws.Rows(2).RowHeight = $origRow2Height
